$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107, shifting existing rows 107-221 down to 108-222
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new data record
$ws.Cells.Item(107, 1).Value = 6
$ws.Cells.Item(107, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(107, 3).Value = "Metropolitana"
$ws.Cells.Item(107, 4).Value = 44484
$ws.Cells.Item(107, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(107, 5).Value = 13
$ws.Cells.Item(107, 6).Value = 100112032
$ws.Cells.Item(107, 7).Value = "Zapallo italiano"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 320
$ws.Cells.Item(107, 11).Value = 13000
$ws.Cells.Item(107, 12).Value = 14000
$ws.Cells.Item(107, 13).Value = 13375
$ws.Cells.Item(107, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(107, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(107, 16).Value = 268
$ws.Cells.Item(107, 17).Value = 50
$ws.Cells.Item(107, 18).Value = "Hortaliza"
